$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update time_taken values (F2:F9) on the "data" sheet ---
$ws.Range("F2").Value = "2021-10-05 14:34:15.808561"
$ws.Range("F3").Value = "2021-10-05 14:34:15.808569"
$ws.Range("F4").Value = "2021-10-05 14:34:15.808572"
$ws.Range("F5").Value = "2021-10-05 14:34:15.808575"
$ws.Range("F6").Value = "2021-10-05 14:34:15.808577"
$ws.Range("F7").Value = "2021-10-05 14:34:15.808580"
$ws.Range("F8").Value = "2021-10-05 14:34:15.808583"
$ws.Range("F9").Value = "2021-10-05 14:34:15.808585"

# --- Add the new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Reuse the same header formatting (bold, centered, bordered) used on the
# "data" sheet's header row / index column, so the same style index is used.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Kabuki syndrome"
$meta.Range("C2").Value = 134

# data_version needs to stay a text value ("0.12") rather than be coerced
# to a number - force text format while writing it, then clear the format
# back so no residual number formatting / quote-prefix style is left on
# the cell.
$dv = $meta.Range("D2")
$dv.NumberFormat = "@"
$dv.Value = "0.12"
$dv.NumberFormat = "General"
$dv.Style = "Normal"

$meta.Range("E2").Value = "2021-03-28T23:23:22.699470Z"
$meta.Range("F2").Value = "2021-10-05 14:34:15.805298"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/134/?format=json"
